$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    'Shawn Mendes',
    'Chris Martin',
    'Santos',
    'Adam Levine',
    'Robbie Williams',
    'Usher',
    'The Weeknd',
    'Taio Cruz',
    'will.i.am',
    'Akon',
    'Bruno Mars',
    'Miguel',
    'Ed Sheeran',
    'Emilio',
    'Justin Bieber',
    'Michael Bublé',
    'Harry Styles',
    'Sam Smith',
    'David Archuleta',
    'Charlie Puth',
    'Taylor Swift',
    'Lana del Rey',
    'Nelly Fortado',
    'Ayliva',
    'Kesha',
    'Charli XCX',
    'Olivia Rodrigo',
    'Katy Perry',
    'Fergie',
    'Caroline Polachek',
    'Christina Aguilera',
    'Adele',
    'Dua Lipa',
    'Lady Gaga',
    'Mariah Carey',
    'Anastacia',
    'Kim Petras',
    'Billie Eilish',
    'Miley Cyrus',
    'Britney Spears',
    'Eminem',
    'Kurt Cobain',
    'Billy Joe Armstrong',
    'Ray Charles',
    'Jim Morrison',
    'Roger Cicero',
    '50 Cent',
    'Keith Bryant',
    'Chris Cagle',
    'Luciano',
    'Josh Gracin',
    'Ludacris',
    'Travis Denning',
    'Geza X',
    'Krizz Kaliko',
    'Benzino',
    'Masta Ace',
    'Cassidy',
    'Don Toliver',
    'Slowthai',
    'Nina Simone',
    'Sarah Vaughan',
    'Diane Chase',
    'Alee',
    'Liz Anderson',
    'Ashley Arrison',
    'Deana Carter',
    'Billie Holiday',
    'Amy Dalley',
    'Roxie Dean',
    'Eva O',
    'Alecia Elliott',
    'Courtney Love',
    'Sierra Kay',
    'Amil',
    'Tiffany Foxx',
    'Mia X',
    'Nonchalant',
    'Khia',
    'Kaliii'
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $names[$i]
}

$ws.Application.ActiveWindow.ScrollRow = 53
$ws.Range("B82").Select()
